$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.005.49"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.243.79"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.65%  "
$ws.Range("E7").Value = "  -1.57%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.542"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.84"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0824"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("E12").Value = "  -2.32%  "
$ws.Range("E13").Value = "  -1.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.588.86"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.845"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.245.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.966.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.25%  "
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("E21").Value = "  -3.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "234.46"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.77%  "
$ws.Range("E25").Value = "  -5.04%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.14%  "
$ws.Range("E28").Value = "  -0.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.78"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.19"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.03"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.28%  "
$ws.Range("E33").Value = "  -3.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.67"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.19"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.113"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.20%  "
$ws.Range("E37").Value = "  +2.34%  "
$ws.Range("E38").Value = "  -2.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.03"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +11.86%  "
$ws.Range("E40").Value = "  -0.63%  "
$ws.Range("E41").Value = "  -5.86%  "
$ws.Range("E42").Value = "  -2.30%  "
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.742.07"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.197"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "81.54"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "74.27"
$ws.Range("D47").Style = "Normal"
$ws.Range("E48").Value = "  -2.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "102.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.65"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "57.42"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.19%  "
